# Update dashboards - 2026-02-06
# Applies the weekly/monthly data refresh to the Aguilar Prototype sheet:
#  - rolls several "as-of" dates (N column) forward to their latest release
#  - shifts the trailing 5-observation history (Q:U) for each series
#  - a few rows only need their highlight style toggled (newly-updated vs
#    no-longer-newest), without any value changes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 10-12 (JOLTS Openings/Hires/Separations Rate): as-of date rolls
# from 2025-11-01 (45962) to 2025-12-01 (45992), and gets the "latest"
# highlight style (48 -> 49). Use N5 (already style 49) as a format
# donor so we reuse the existing style index instead of minting a new one.
$ws.Range("N5").Copy()
$ws.Range("N10").PasteSpecial(-4122)
$ws.Range("N11").PasteSpecial(-4122)
$ws.Range("N12").PasteSpecial(-4122)

$ws.Range("N10").Value = 45992
$ws.Range("Q10").Value = 3.9
$ws.Range("R10").Value = 4.2
$ws.Range("S10").Value = 4.5
$ws.Range("T10").Value = 4.6
$ws.Range("U10").Value = 4.3

$ws.Range("N11").Value = 45992
$ws.Range("Q11").Value = 3.3
$ws.Range("R11").Value = 3.2
$ws.Range("S11").Value = 3.4
$ws.Range("T11").Value = 3.4
$ws.Range("U11").Value = 3.2

$ws.Range("N12").Value = 45992
$ws.Range("Q12").Value = 3.3
$ws.Range("R12").Value = 3.2
$ws.Range("S12").Value = 3.2
$ws.Range("T12").Value = 3.3
$ws.Range("U12").Value = 3.2

# ---------------------------------------------------------------------
# Rows 22-23 (PPI-FD): same as-of date (45992), but these lose the
# "latest" highlight (49 -> 48). N3 is already plain style 48.
$ws.Range("N3").Copy()
$ws.Range("N22").PasteSpecial(-4122)
$ws.Range("N23").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Rows 29-30 (5yr5yr Forward / 10yr TIPS): as-of date rolls one day
# (46057 -> 46058); style unchanged.
$ws.Range("N29").Value = 46058
$ws.Range("Q29").Value = 2.16
$ws.Range("R29").Value = 2.19
$ws.Range("S29").Value = 2.19
$ws.Range("T29").Value = 2.18
$ws.Range("U29").Value = 2.19

$ws.Range("N30").Value = 46058
$ws.Range("Q30").Value = 2.32
$ws.Range("R30").Value = 2.35
$ws.Range("S30").Value = 2.36
$ws.Range("T30").Value = 2.35
$ws.Range("U30").Value = 2.36

# ---------------------------------------------------------------------
# Rows 47-50, 52 (FFR, 2y/5y/10y UST, BAA): as-of date rolls one day
# (46056 -> 46057); style unchanged.
$ws.Range("N47").Value = 46057

$ws.Range("N48").Value = 46057
$ws.Range("S48").Value = 3.57
$ws.Range("T48").Value = 3.52
$ws.Range("U48").Value = 3.53

$ws.Range("N49").Value = 46057
$ws.Range("S49").Value = 3.83
$ws.Range("T49").Value = 3.79
$ws.Range("U49").Value = 3.8

$ws.Range("N50").Value = 46057
$ws.Range("Q50").Value = 4.29
$ws.Range("R50").Value = 4.28
$ws.Range("S50").Value = 4.29
$ws.Range("T50").Value = 4.26
$ws.Range("U50").Value = 4.24

$ws.Range("N52").Value = 46057
$ws.Range("Q52").Value = 5.93
$ws.Range("R52").Value = 5.91
$ws.Range("S52").Value = 5.9
$ws.Range("T52").Value = 5.88
$ws.Range("U52").Value = 5.87

# ---------------------------------------------------------------------
# Row 51 (30y Mortgage): as-of date rolls from 2026-01-26 (46048) to
# 2026-02-02 (46055), and gets the "latest" highlight (48 -> 49).
$ws.Range("N5").Copy()
$ws.Range("N51").PasteSpecial(-4122)

$ws.Range("N51").Value = 46055
$ws.Range("Q51").Value = 6.11
$ws.Range("R51").Value = 6.1
$ws.Range("S51").Value = 6.09
$ws.Range("T51").Value = 6.06
$ws.Range("U51").Value = 6.16
